# Update phyFa OS code size section: add a new "Phyfa.OS" block
# (columns E/F) alongside the existing FreeRTOS block in columns B/C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new block, mirroring C1's "FreeRTOS" header but in column F.
$ws.Range("F1").Value = "Phyfa.OS"

# Row 5: OS lib code size label + value.
$ws.Range("E5").Value = "OS.lib code size"
$ws.Range("F5").Value = "11,992Bytes"

# Row 6: OAL(OSA) code size label + value.
$ws.Range("E6").Value = "OAL(OSA)code size"
$ws.Range("F6").Value = "30,032Bytes"

# Give the new columns sensible widths, matching the sizing used for the
# existing table (D is a spacer column, E/F hold the new label/value pair).
$ws.Columns.Item(4).ColumnWidth = 23.285714285714285
$ws.Columns.Item(5).ColumnWidth = 20.504464285714285
$ws.Columns.Item(6).ColumnWidth = 14.840401785714286

# Reflect the new selection left behind by the edit.
$ws.Range("C14").Select()
